# Adds a "2022-Q1" sheet (fund-holding detail) before the "总计" sheet,
# and records the new quarter's summary row on "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value as TEXT (not auto-coerced to a number) without
# picking up stray number-formatting. Excel treats a leading apostrophe
# as a "force text" marker and strips it from the stored value.
# ---------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = "'" + $val
}

# =======================================================================
# 1. Insert the new "2022-Q1" worksheet immediately before "总计"
#
# NOTE: worksheet variables in this host resolve by *position*, not
# stable identity - holding a reference across a Worksheets.Add() call
# (which shifts indices) silently repoints it at whatever sheet now
# occupies that slot. So: do the add/rename in total isolation, then
# re-fetch every handle we need by name afterwards, once the sheet
# collection is done changing shape.
# =======================================================================
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Collection shape is now final - safe to grab stable handles by name.
$newSheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")
$styleSource = $wb.Worksheets.Item(1)   # "2020-Q4" - same table layout/style to copy from

# =======================================================================
# 2. Fund holding rows for 2022-Q1
# =======================================================================
$fundRows = @(
    ,@(0, "007119", "睿远成长价值混合A", "280.90", "92.63", "3.52", "9.8877", 10)
    ,@(1, "011164", "富国兴远优选12个月持有期混合A", "84.49", "82.18", "9.05", "7.6463", 1)
    ,@(2, "001186", "富国文体健康股票A", "18.44", "86.61", "9.55", "1.7610", 1)
    ,@(3, "011165", "富国兴远优选12个月持有期混合C", "18.89", "82.18", "9.05", "1.7095", 1)
    ,@(4, "007120", "睿远成长价值混合C", "29.98", "92.63", "3.52", "1.0553", 10)
    ,@(5, "000586", "景顺长城中小板创业板精选股票", "2.42", "94.15", "6.19", "0.1498", 6)
    ,@(6, "009606", "长信稳健精选混合A", "7.44", "37.23", "1.33", "0.0990", 6)
    ,@(7, "011125", "富国文体健康股票C", "0.71", "86.61", "9.55", "0.0678", 1)
    ,@(8, "010706", "景顺长城景骊成长混合型证券投资基金", "1.13", "93.50", "5.62", "0.0635", 8)
    ,@(9, "519963", "长信利盈灵活配置混合A", "4.82", "29.86", "1.15", "0.0554", 9)
    ,@(10, "519961", "长信利广灵活配置混合A", "4.50", "41.59", "1.15", "0.0518", 8)
    ,@(11, "519962", "长信利盈灵活配置混合C", "3.02", "29.86", "1.15", "0.0347", 9)
    ,@(12, "161224", "国投瑞银新丝路灵活配置混合(LOF)", "0.77", "94.48", "4.40", "0.0339", 6)
    ,@(13, "009607", "长信稳健精选混合C", "1.59", "37.23", "1.33", "0.0211", 6)
    ,@(14, "005444", "光大保德信多策略精选18个月定期开放灵活配置混合", "1.00", "29.09", "1.67", "0.0167", 10)
    ,@(15, "519960", "长信利广灵活配置混合C", "0.70", "41.59", "1.15", "0.0080", 8)
    ,@(16, "006346", "安信量化优选股票A", "0.71", "90.62", "0.67", "0.0048", 5)
    ,@(17, "006347", "安信量化优选股票C", "0.49", "90.62", "0.67", "0.0033", 5)
    ,@(18, "010999", "兴华瑞丰混合A", "0.06", "29.21", "2.82", "0.0017", 6)
    ,@(19, "011000", "兴华瑞丰混合C", "0.05", "29.21", "2.82", "0.0014", 6)
)

# Header row (B1:H1)
Set-TextCell $newSheet 1 2 "基金代码"
Set-TextCell $newSheet 1 3 "基金名称"
Set-TextCell $newSheet 1 4 "基金规模"
Set-TextCell $newSheet 1 5 "股票总仓位"
Set-TextCell $newSheet 1 6 "仓位占比"
Set-TextCell $newSheet 1 7 "持有市值(亿元)"
Set-TextCell $newSheet 1 8 "仓位排名"

# Data rows
$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $newSheet $r 2 $row[1]
    Set-TextCell $newSheet $r 3 $row[2]
    Set-TextCell $newSheet $r 4 $row[3]
    Set-TextCell $newSheet $r 5 $row[4]
    Set-TextCell $newSheet $r 6 $row[5]
    Set-TextCell $newSheet $r 7 $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Match the bold/centered/bordered header style used on the other
# quarterly sheets (copy the format only, values are untouched).
$styleSource.Cells.Item(1, 2).Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Cells.Item(2, 1).Copy()
$newSheet.Range("A2:A21").PasteSpecial(-4122)

# =======================================================================
# 3. Add the 2022-Q1 summary row to "总计" (new row 2, pushing the rest down)
# =======================================================================
$totalSheet.Rows.Item(2).Insert()

# Re-apply the bold/centered/bordered style to the new A2 (index column)
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)
# Plain data style (no border/bold) for the rest of the new row
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Cells.Item(2, 1).Value = 0
Set-TextCell $totalSheet 2 2 "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 20
$totalSheet.Cells.Item(2, 4).Value = 22.67
